$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

function Set-Text($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 2
Set-PlainText "D2" ("26.781.61")
Set-Text "E2" ("  -1.43%  ")

# Row 3
Set-PlainText "D3" ("1.547.11")
Set-Text "E3" ("  -1.66%  ")

# Row 4
Set-Text "E4" ("  -0.02%  ")

# Row 5
Set-PlainText "D5" ("205.89")
Set-Text "E5" ("  -0.70%  ")

# Row 6
Set-PlainText "D6" ("0.481")
Set-Text "E6" ("  -1.76%  ")

# Row 7
Set-Text "E7" ("  +0.01%  ")

# Row 8
Set-Text "B8" ("Solana")
Set-Text "C8" ("https://coinranking.com/coin/zNZHO_Sjf+solana-sol")
Set-PlainText "D8" ("21.40")
Set-Text "E8" ("  -3.91%  ")

# Row 9
Set-Text "B9" ("Cardano")
Set-Text "C9" ("https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada")
Set-PlainText "D9" ("0.246")
Set-Text "E9" ("  -1.15%  ")

# Row 10
Set-PlainText "D10" ("0.0581")
Set-Text "E10" ("  -1.36%  ")

# Row 11
Set-PlainText "D11" ("0.0853")
Set-Text "E11" ("  -1.72%  ")

# Row 12
Set-PlainText "D12" ("1.770.33")
Set-Text "E12" ("  -1.45%  ")

# Row 13
Set-PlainText "D13" ("1.554.65")
Set-Text "E13" ("  -1.16%  ")

# Row 14
Set-PlainText "D14" ("3.67")
Set-Text "E14" ("  -2.70%  ")

# Row 15
Set-PlainText "D15" ("0.511")
Set-Text "E15" ("  -1.28%  ")

# Row 16
Set-PlainText "D16" ("26.790.84")
Set-Text "E16" ("  -1.46%  ")

# Row 17
Set-PlainText "D17" ("61.12")
Set-Text "E17" ("  -1.82%  ")

# Row 18
Set-PlainText "D18" ("213.81")
Set-Text "E18" ("  -0.23%  ")

# Row 19
Set-PlainText "D19" ("0.0" + ([string][char]0x2083) + "0687")
Set-Text "E19" ("  +0.26%  ")

# Row 20
Set-PlainText "D20" ("7.24")
Set-Text "E20" ("  -1.52%  ")

# Row 21
Set-Text "E21" ("  +0.01%  ")

# Row 22
Set-PlainText "D22" ("4.08")
Set-Text "E22" ("  -0.91%  ")

# Row 23
Set-PlainText "D23" ("8.98")
Set-Text "E23" ("  -4.35%  ")

# Row 24
Set-PlainText "D24" ("1.98")
Set-Text "E24" ("  -1.36%  ")

# Row 25
Set-PlainText "D25" ("153.42")
Set-Text "E25" ("  +0.57%  ")

# Row 26
Set-PlainText "D26" ("6.51")
Set-Text "E26" ("  -2.54%  ")

# Row 27
Set-PlainText "D27" ("14.90")
Set-Text "E27" ("  -0.28%  ")

# Row 28
Set-Text "E28" ("  -0.06%  ")

# Row 29
Set-Text "E29" ("  -1.46%  ")

# Row 30
Set-Text "E30" ("  -0.53%  ")

# Row 31
Set-Text "E31" ("  -0.90%  ")

# Row 32
Set-PlainText "D32" ("3.18")
Set-Text "E32" ("  +0.08%  ")

# Row 33
Set-PlainText "D33" ("1.345.16")
Set-Text "E33" ("  -4.60%  ")

# Row 34
Set-Text "E34" ("  -0.27%  ")

# Row 35
Set-Text "E35" ("  -3.41%  ")

# Row 36
Set-Text "E36" ("  -0.42%  ")

# Row 37
Set-PlainText "D37" ("0.928")
Set-Text "E37" ("  -1.28%  ")

# Row 38
Set-Text "E38" ("  -1.06%  ")

# Row 39
Set-PlainText "D39" ("0.519")
Set-Text "E39" ("  +0.67%  ")

# Row 40
Set-PlainText "D40" ("0.800")
Set-Text "E40" ("  -1.84%  ")

# Row 41
Set-Text "B41" ("FraxShare")
Set-Text "C41" ("https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs")
Set-PlainText "D41" ("5.67")
Set-Text "E41" ("  +5.11%  ")

# Row 42
Set-Text "B42" ("WEMIXToken")
Set-Text "C42" ("https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix")
Set-PlainText "D42" ("0.993")
Set-Text "E42" ("  -0.47%  ")

# Row 43
Set-Text "B43" ("MXToken")
Set-Text "C43" ("https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx")
Set-PlainText "D43" ("2.19")
Set-Text "E43" ("  -0.45%  ")

# Row 44
Set-Text "B44" ("RenderToken")
Set-Text "C44" ("https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr")
Set-PlainText "D44" ("1.76")
Set-Text "E44" ("  -3.92%  ")

# Row 45
Set-Text "B45" ("Aave")
Set-Text "C45" ("https://coinranking.com/coin/ixgUfzmLR+aave-aave")
Set-PlainText "D45" ("62.84")
Set-Text "E45" ("  -1.36%  ")

# Row 46
Set-Text "B46" ("RocketPoolETH")
Set-Text "C46" ("https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth")
Set-PlainText "D46" ("1.684.17")
Set-Text "E46" ("  -1.46%  ")

# Row 47
Set-Text "B47" ("mCoin")
Set-Text "C47" ("https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin")
Set-PlainText "D47" ("2.25")
Set-Text "E47" ("  -2.94%  ")

# Row 48
Set-Text "B48" ("Quant")
Set-Text "C48" ("https://coinranking.com/coin/bauj_21eYVwso+quant-qnt")
Set-PlainText "D48" ("85.73")
Set-Text "E48" ("  +0.01%  ")

# Row 49
Set-Text "B49" ("Cronos")
Set-Text "C49" ("https://coinranking.com/coin/65PHZTpmE55b+cronos-cro")
Set-PlainText "D49" ("0.0512")
Set-Text "E49" ("  +3.60%  ")

# Row 50
Set-Text "B50" ("BabyDogeCoin")
Set-Text "C50" ("https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge")
Set-PlainText "D50" ("0.0" + ([string][char]0x2087) + "0972")
Set-Text "E50" ("  -0.94%  ")

# Row 51
Set-Text "B51" ("Algorand")
Set-Text "C51" ("https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo")
Set-PlainText "D51" ("0.0949")
Set-Text "E51" ("  -0.45%  ")
